# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets, reflecting freshly generated data from the gh-pages build at 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new value mapping for the "展览" sheet (column F)
$exhibitUpdates = @{
    4  = 579
    5  = 1824
    9  = 2374
    12 = 156
    13 = 1436
    14 = 507
    18 = 18
    20 = 197
    24 = 94
    26 = 1483
    28 = 372
    29 = 234
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new value mapping for the "全部类型" sheet (column F)
$allTypesUpdates = @{
    4  = 579
    5  = 1824
    10 = 2374
    13 = 156
    14 = 1436
    15 = 507
    19 = 18
    21 = 197
    25 = 94
    27 = 1483
    29 = 372
    30 = 234
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}
